$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.903.65"
$ws.Range("E2").Value = "  -0.23%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.862.90"
$ws.Range("E3").Value = "  +0.14%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.91"
$ws.Range("E5").Value = "  -0.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("E6").Value = "  -0.07%  "

$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3649"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07170"
$ws.Range("E9").Value = "  +0.52%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8931"
$ws.Range("E10").Value = "  +0.62%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.69"
$ws.Range("E11").Value = "  +0.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.879.49"
$ws.Range("E12").Value = "  +1.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07500"
$ws.Range("E13").Value = "  -0.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "95.13"
$ws.Range("E14").Value = "  +6.88%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.235"
$ws.Range("E15").Value = "  -1.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008516"
$ws.Range("E17").Value = "  +1.82%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.22"
$ws.Range("E18").Value = "  +1.12%  "

$ws.Range("E19").Value = "  -0.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.944.32"
$ws.Range("E20").Value = "  -0.23%  "

$ws.Range("E21").Value = "  -0.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.117.46"
$ws.Range("E22").Value = "  +1.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.40"
$ws.Range("E23").Value = "  -0.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.409"
$ws.Range("E24").Value = "  -0.79%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.97"
$ws.Range("E25").Value = "  +0.57%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.780"
$ws.Range("E26").Value = "  -3.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.89"
$ws.Range("E27").Value = "  -0.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.088"
$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.24"
$ws.Range("E29").Value = "  +0.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.710"
$ws.Range("E30").Value = "  +1.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.669"
$ws.Range("E31").Value = "  +0.67%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09225"
$ws.Range("E32").Value = "  +2.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05145"
$ws.Range("E33").Value = "  +0.71%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7511"
$ws.Range("E34").Value = "  +3.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.962"
$ws.Range("E35").Value = "  -3.03%  "

$ws.Range("E36").Value = "  +0.44%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.257"
$ws.Range("E37").Value = "  +7.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.588"
$ws.Range("E38").Value = "  +5.62%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02002"
$ws.Range("E39").Value = "  -1.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5574"
$ws.Range("E40").Value = "  +4.86%  "

$ws.Range("E41").Value = "  -0.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.574"
$ws.Range("E42").Value = "  -0.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.68"
$ws.Range("E43").Value = "  +1.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.616"
$ws.Range("E44").Value = "  +4.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1472"
$ws.Range("E45").Value = "  +0.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4696"
$ws.Range("E46").Value = "  +1.97%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9997"
$ws.Range("E47").Value = "  -0.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.03"
$ws.Range("E48").Value = "  +0.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.561"
$ws.Range("E49").Value = "  +0.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.71"
$ws.Range("E50").Value = "  +0.58%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.06"
$ws.Range("E51").Value = "  -1.34%  "
